$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2021_13")
$ws.Range("C13").Value = "7.9660"
$ws.Range("E13").Value = "79.660"
$ws.Range("C25").Value = "3.1674"
$ws.Range("E25").Value = "31.674"
$ws.Range("C37").Value = "0.3976"
$ws.Range("E37").Value = "3.976"

$ws = $wb.Worksheets.Item("2021_24")
$ws.Range("C4").Value = "1.0814"
$ws.Range("D4").Value = "1.052"
$ws.Range("E4").Value = "1.112"
$ws.Range("C13").Value = "0.4953"
$ws.Range("D13").Value = "0.196"
$ws.Range("E13").Value = "1.250"
$ws.Range("C16").Value = "1.1986"
$ws.Range("D16").Value = "1.173"
$ws.Range("E16").Value = "1.225"
$ws.Range("C25").Value = "1.1109"
$ws.Range("D25").Value = "0.450"
$ws.Range("E25").Value = "2.740"
$ws.Range("C28").Value = "1.1084"
$ws.Range("D28").Value = "1.077"
$ws.Range("E28").Value = "1.141"
$ws.Range("C37").Value = "2.2428"
$ws.Range("D37").Value = "0.638"
$ws.Range("E37").Value = "7.890"

$ws = $wb.Worksheets.Item("2022_06")
$ws.Range("C4").Value = "1.0350"
$ws.Range("C13").Value = "2.0390"
$ws.Range("D13").Value = "0.830"
$ws.Range("E13").Value = "5.010"
$ws.Range("C16").Value = "1.0718"
$ws.Range("D16").Value = "1.048"
$ws.Range("E16").Value = "1.096"
$ws.Range("C25").Value = "1.6744"
$ws.Range("D25").Value = "1.140"
$ws.Range("E25").Value = "2.459"
$ws.Range("C28").Value = "1.0355"
$ws.Range("E28").Value = "1.083"
$ws.Range("C37").Value = "0.8212"
$ws.Range("D37").Value = "0.326"
$ws.Range("E37").Value = "2.067"
$ws.Range("C40").Value = "1.6281"
$ws.Range("D40").Value = "1.596"
$ws.Range("E40").Value = "1.661"
$ws.Range("C49").Value = "3.4111"
$ws.Range("D49").Value = "2.008"
$ws.Range("E49").Value = "5.795"
$ws.Range("C52").Value = "1.5190"
$ws.Range("D52").Value = "1.488"
$ws.Range("E52").Value = "1.551"
$ws.Range("C61").Value = "2.0372"
$ws.Range("D61").Value = "1.152"
$ws.Range("E61").Value = "3.603"

$ws = $wb.Worksheets.Item("2022_47")
$ws.Range("C13").Value = "4.1845"
$ws.Range("E13").Value = "41.845"
$ws.Range("C25").Value = "4.9146"
$ws.Range("E25").Value = "49.146"
$ws.Range("C37").Value = "3.5630"
$ws.Range("E37").Value = "35.630"
$ws.Range("C49").Value = "3.7594"
$ws.Range("E49").Value = "37.594"
